$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple B-column updates (Taxonsorteringsordning) ---
$ws.Range("B2").Value = 79243
$ws.Range("B3").Value = 79243
$ws.Range("B4").Value = 57881
$ws.Range("B7").Value = 79243
$ws.Range("B8").Value = 79243
$ws.Range("B9").Value = 79243
$ws.Range("B10").Value = 79243
$ws.Range("B11").Value = 79243
$ws.Range("B12").Value = 79243
$ws.Range("B13").Value = 57073
$ws.Range("B14").Value = 79243
$ws.Range("B15").Value = 57881
$ws.Range("B16").Value = 79243
$ws.Range("B17").Value = 79243

# --- Rows 5 and 6 effectively swap their species data, with new B values ---
$ws.Range("A5").Value = 130981914
$ws.Range("B5").Value = 57881
$ws.Range("E5").Value = 100049
$ws.Range("F5").Value = "Spillkråka"
$ws.Range("G5").Value = "Dryocopus martius"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("J5").ClearContents()
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "äldre spår"
$ws.Range("Q5").Value = 437688
$ws.Range("R5").Value = 6792409
$ws.Range("AF5").ClearContents()
$ws.Range("AX5").Value = "Eva Löfqvist, Alfhild Sehlin"

$ws.Range("A6").Value = 130981911
$ws.Range("B6").Value = 91829
$ws.Range("E6").Value = 5442
$ws.Range("F6").Value = "Tallticka"
$ws.Range("G6").Value = "Porodaedalea pini"
$ws.Range("H6").Value = "(Brot.) Murrill"
$ws.Range("J6").Value = ""
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value = 437697
$ws.Range("R6").Value = 6792416
$ws.Range("AF6").Value = ""
$ws.Range("AX6").Value = "Eva Löfqvist"
